$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "앱코 A660 3325 프로페셔널 게이밍 마우스"
$ws.Range("B3").Value = "ROCCAT KONE PURE SEL RGB 게이밍 마우스"
$ws.Range("B5").Value = "COX CM50 RGB 게이밍 마우스"
$ws.Range("B7").Value = "앱코 A200 3050 RGB 게이밍 마우스"
$ws.Range("B10").Value = "TG삼보 TM615U HEALING 인체공학 버티컬 유선마우스"
$ws.Range("B11").Value = "앱코 A700 하이엔드 게이밍 마우스"
$ws.Range("B13").Value = "로지텍G G PRO HERO 게이밍 마우스"
$ws.Range("B14").Value = "앱코 A250 3050 RGB 게이밍 마우스"
$ws.Range("B15").Value = "앱코 AEM20 인체공학 버티컬 유선 마우스"
$ws.Range("B17").Value = "앱코 A660 3360 하이엔드 게이밍 마우스"
$ws.Range("B19").Value = "COX CM60 TARGET RGB 게이밍 마우스"
$ws.Range("B20").Value = "플레오맥스 MO-ER50 버티컬 유선 마우스"
$ws.Range("B22").Value = "리줌 M1 RGB 매크로 버티컬 마우스"
$ws.Range("B23").Value = "TG삼보 TG-M300UN 인체공학 무소음 유선 마우스"
$ws.Range("B25").Value = "블레스정보통신 지오 i980 RGB 버티컬 인체공학 마우스"
$ws.Range("B27").Value = "COX CM90 RGB 게이밍 마우스"
$ws.Range("B28").Value = "TG삼보 TM137U HEALING 인체공학 버티컬마우스"
$ws.Range("B29").Value = "앱코 A900 3389 RGB 게이밍 마우스"
$ws.Range("B30").Value = "세컨드찬스 긱스타 GM900 3325 LED 게이밍 마우스"
$ws.Range("B31").Value = "HP M100S Gaming Mouse"
$ws.Range("B32").Value = "다얼유 LM121 RGB 버티컬 인체공학 마우스"
$ws.Range("B34").Value = "앱코 A800 3389 초경량 RGB 게이밍 마우스"
$ws.Range("B35").Value = "COX CM10 RGB 게이밍 마우스"
$ws.Range("B36").Value = "리줌 RS-G1 RGB 게이밍 마우스"
$ws.Range("B37").Value = "필립스 M444 유선 버티컬 마우스"
$ws.Range("B39").Value = "앱코 A100 레인보우 LED 게이밍 마우스"
$ws.Range("B41").Value = "마이크로소프트 프로 인텔리 마우스"
$ws.Range("B42").Value = "로지텍 B100 Optical Mouse"
$ws.Range("B45").Value = "펄사 Xlite 초경량 유선 게이밍  마우스"
$ws.Range("B46").Value = "HP M200 Gaming Mouse"
$ws.Range("B47").Value = "청연엠엔에스 NAVEE NV73-VMS10U 버티컬 무소음 유선마우스"
$ws.Range("B49").Value = "펄사 Xlite 초경량 유선 게이밍 마우스 풀세트"
$ws.Range("B50").Value = "마이크로닉스 MANIC G40 RGB PMW3360 게이밍 마우스"
$ws.Range("B52").Value = "로지텍 G PRO 게이밍 마우스"
$ws.Range("B53").Value = "삼성전자 삼성 SPA-MMG1PUB 게이밍마우스"
$ws.Range("B55").Value = "HP M280 Gaming Mouse"
$ws.Range("B60").Value = "플레오맥스 MO-ER700 인체공학 버티칼 마우스"
$ws.Range("B61").Value = "제닉스 STORMX M1 게이밍 마우스"
$ws.Range("B63").Value = "다얼유 EM925 RGB 게이밍 마우스 다이아몬드 에디션"
$ws.Range("B64").Value = "맥스틸 TRON G20 PRO PMW 3330 RGB 게이밍 마우스"
$ws.Range("B67").Value = "앱코 H1 쉘체인저 게이밍 마우스"
$ws.Range("B69").Value = "COSY 버티컬 이지 마우스(M3124)"
$ws.Range("B71").Value = "앱코 A810 3327 RGB 피어스 게이밍 마우스"
$ws.Range("B73").Value = "앱코 MX1200 옵티컬 게이밍 마우스"
$ws.Range("B75").Value = "알텍랜싱 ALBM7444 유선 버티컬 마우스"
$ws.Range("B76").Value = "HP M270 Gaming Mouse"
$ws.Range("B77").Value = "리줌 M8 유선 버티컬 마우스"
$ws.Range("B78").Value = "맥스틸 RATIO S10 PMW 3389 게이밍 마우스"
$ws.Range("B79").Value = "마이크로소프트 에고노믹 마우스"
$ws.Range("B80").Value = "COSY 버티컬 피트 마우스(M1189)"
$ws.Range("B85").Value = "ASUS ROG GLADIUS II ORIGIN PINK 게이밍 마우스"
$ws.Range("B87").Value = "지클릭커 GM-M250 LED 무소음 게이밍 마우스"
$ws.Range("B88").Value = "제닉스 TITAN G 게이밍 마우스"
$ws.Range("B91").Value = "HP M260 Gaming Mouse"
$ws.Range("B92").Value = "TG삼보 TG-M300U 인체공학 유선 마우스"
$ws.Range("B93").Value = "앱코 HACKER A530 3325 RGB 게이밍 마우스"
$ws.Range("B94").Value = "HP G360 Gaming Mouse"
$ws.Range("B101").Value = "COX CM80 3330 RGB 게이밍 마우스"
